$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 59950
$ws.Range("J32").Value = 59950
$ws.Range("L32").Value = 59950
$ws.Range("N32").Value = -60602
$ws.Range("H98").Value = 583.7857
$ws.Range("I98").Value = 583.7857
$ws.Range("K98").Value = 583.7857
$ws.Range("M98").Value = 914.2143
$ws.Range("H108").Value = 99278
$ws.Range("J108").Value = 99278
$ws.Range("L108").Value = 99278
$ws.Range("N108").Value = -106958
$ws.Range("H110").Value = 57704.855
$ws.Range("J110").Value = 57704.855
$ws.Range("L110").Value = 57704.855
$ws.Range("N110").Value = -65884.85500000001
$ws.Range("H117").Value = 90738.8
$ws.Range("J117").Value = 90738.8
$ws.Range("L117").Value = 90738.8
$ws.Range("N117").Value = -99916.8
$ws.Range("H120").Value = 49189.6
$ws.Range("J120").Value = 49189.6
$ws.Range("L120").Value = 49189.6
$ws.Range("N120").Value = -58865.6
$ws.Range("H122").Value = 583.7857
$ws.Range("I122").Value = 583.7857
$ws.Range("K122").Value = 1751.3571
$ws.Range("M122").Value = 698.6428999999998
$ws.Range("H123").Value = 61712
$ws.Range("J123").Value = 60762.855
$ws.Range("L123").Value = 60762.855
$ws.Range("N123").Value = -70562.85500000001
$ws.Range("H132").Value = 2089.25
$ws.Range("I132").Value = 1931.5
$ws.Range("K132").Value = 5794.5
$ws.Range("M132").Value = -3264.5
$ws.Range("H133").Value = 76252.63
$ws.Range("J133").Value = 76252.63
$ws.Range("L133").Value = 76252.63
$ws.Range("N133").Value = -86372.63
$ws.Range("H134").Value = 98962.86
$ws.Range("J134").Value = 98962.86
$ws.Range("L134").Value = 98962.86
$ws.Range("N134").Value = -109102.86
$ws.Range("H137").Value = 694290.1
$ws.Range("I137").Value = 3371.3333
$ws.Range("K137").Value = 10113.9999
$ws.Range("M137").Value = -7563.999899999999
$ws.Range("H138").Value = 2088.25
$ws.Range("J138").Value = 2499.25
$ws.Range("L138").Value = 7497.75
$ws.Range("N138").Value = -17777.75
$ws.Range("H139").Value = 74505.22
$ws.Range("J139").Value = 74505.22
$ws.Range("L139").Value = 74505.22
$ws.Range("N139").Value = -84785.22
$ws.Range("H140").Value = 66637.875
$ws.Range("J140").Value = 71982.336
$ws.Range("L140").Value = 71982.336
$ws.Range("N140").Value = -82342.336
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 286.125
$ws.Range("I5").Value = 304.7143
$ws.Range("K5").Value = 304.7143
$ws.Range("M5").Value = -192.7143
$ws.Range("H32").Value = 4994.0684
$ws.Range("I32").Value = 1805.6666
$ws.Range("K32").Value = 1805.6666
$ws.Range("M32").Value = -1518.6666
$ws.Range("H45").Value = 41671668
$ws.Range("I45").Value = 5000
$ws.Range("K45").Value = 5000
$ws.Range("M45").Value = -4623
$ws.Range("H52").Value = 49712.668
$ws.Range("J52").Value = 49712.668
$ws.Range("L52").Value = 49712.668
$ws.Range("N52").Value = -50348.668
$ws.Range("H97").Value = 680.17645
$ws.Range("I97").Value = 703.875
$ws.Range("J97").Value = 301
$ws.Range("K97").Value = 703.875
$ws.Range("L97").Value = 301
$ws.Range("M97").Value = -207.875
$ws.Range("N97").Value = -1293
$ws.Range("H119").Value = 41599.8
$ws.Range("J119").Value = 41599.8
$ws.Range("L119").Value = 41599.8
$ws.Range("N119").Value = -51275.8
$ws.Range("H121").Value = 86060.25
$ws.Range("J121").Value = 86060.25
$ws.Range("L121").Value = 86060.25
$ws.Range("N121").Value = -89554.25
$ws.Range("H125").Value = 86666.336
$ws.Range("J125").Value = 86666.336
$ws.Range("L125").Value = 86666.336
$ws.Range("N125").Value = -96506.336
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 286.125
$ws.Range("I4").Value = 304.7143
$ws.Range("K4").Value = 304.7143
$ws.Range("M4").Value = -189.7143
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").Value = ""
$ws.Range("H51").Value = 47858.4
$ws.Range("J51").Value = 47858.4
$ws.Range("L51").Value = 47858.4
$ws.Range("N51").Value = -48840.4
$ws.Range("H53").Value = 31398
$ws.Range("J53").Value = 31398
$ws.Range("L53").Value = 31398
$ws.Range("N53").Value = -32546
$ws.Range("H94").Value = 3739.6667
$ws.Range("I94").Value = 2957.125
$ws.Range("K94").Value = 2957.125
$ws.Range("M94").Value = -2506.125
$ws.Range("H110").Value = 99985.664
$ws.Range("J110").Value = 99985.664
$ws.Range("L110").Value = 99985.664
$ws.Range("N110").Value = -108165.664
$ws.Range("H114").Value = 99962.336
$ws.Range("J114").Value = 99962.336
$ws.Range("L114").Value = 99962.336
$ws.Range("N114").Value = -108640.336
$ws.Range("H116").Value = 77326.664
$ws.Range("J116").Value = 77326.664
$ws.Range("L116").Value = 77326.664
$ws.Range("N116").Value = -86504.664
$ws.Range("H117").Value = 99961.42999999999
$ws.Range("J117").Value = 99961.42999999999
$ws.Range("L117").Value = 99961.42999999999
$ws.Range("N117").Value = -109139.43
$ws.Range("H118").Value = 60974.11
$ws.Range("J118").Value = 61652.625
$ws.Range("L118").Value = 61652.625
$ws.Range("N118").Value = -64966.625
$ws.Range("H122").Value = 59121.875
$ws.Range("J122").Value = 59121.875
$ws.Range("L122").Value = 59121.875
$ws.Range("N122").Value = -68921.875
$ws.Range("H132").Value = 29817.818
$ws.Range("J132").Value = 29817.818
$ws.Range("L132").Value = 29817.818
$ws.Range("N132").Value = -39937.818
$ws.Range("H134").Value = 4135.6895
$ws.Range("I134").Value = 2749.7144
$ws.Range("K134").Value = 8249.143199999999
$ws.Range("M134").Value = -5714.143199999999
$ws.Range("H140").Value = 43500
$ws.Range("J140").Value = 43500
$ws.Range("L140").Value = 43500
$ws.Range("N140").Value = -53860
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H116").Value = 70603.2
$ws.Range("J116").Value = 70603.2
$ws.Range("L116").Value = 70603.2
$ws.Range("N116").Value = -79781.2
$ws.Range("H117").Value = 32956.152
$ws.Range("J117").Value = 32956.152
$ws.Range("L117").Value = 32956.152
$ws.Range("N117").Value = -42134.152
$ws.Range("H138").Value = 53675
$ws.Range("J138").Value = 54900
$ws.Range("L138").Value = 54900
$ws.Range("N138").Value = -65180
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 7335.3335
$ws.Range("I44").Value = 7335.3335
$ws.Range("K44").Value = 22006.0005
$ws.Range("M44").Value = -21608.0005
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6937.5
$ws.Range("I80").Value = 6937.5
$ws.Range("K80").Value = 6937.5
$ws.Range("M80").Value = -5939.5
$ws.Range("H83").Value = 6937.5
$ws.Range("I83").Value = 6937.5
$ws.Range("K83").Value = 34687.5
$ws.Range("M83").Value = -29695.5
$ws.Range("H109").Value = 24056.6
$ws.Range("J109").Value = 24056.6
$ws.Range("L109").Value = 24056.6
$ws.Range("N109").Value = -26136.6
$ws.Range("H116").Value = 56340.4
$ws.Range("J116").Value = 58886.75
$ws.Range("L116").Value = 58886.75
$ws.Range("N116").Value = -68064.75
$ws.Range("H135").Value = 25000
$ws.Range("J135").Value = 25000
$ws.Range("L135").Value = 25000
$ws.Range("N135").Value = -35140
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H123").Value = 72151.664
$ws.Range("J123").Value = 74871.875
$ws.Range("L123").Value = 74871.875
$ws.Range("N123").Value = -84671.875
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H121").Value = 98996.664
$ws.Range("J121").Value = 98996.664
$ws.Range("L121").Value = 98996.664
$ws.Range("N121").Value = -102490.664
$ws.Range("H127").Value = 85194.5
$ws.Range("J127").Value = 109999
$ws.Range("L127").Value = 109999
$ws.Range("N127").Value = -119919

Write-Host "Applied all cell updates"